$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reassign the "Expert" column (J) reviewer names for rows 3, 5, 7, 9
$ws.Range("J3").Value = "Sharon Wu"
$ws.Range("J5").Value = "Zhi Zhang"
$ws.Range("J7").Value = "Bo Cao"
$ws.Range("J9").Value = "Peter Michalski"

# Update the active selection to reflect where the user last clicked
$ws.Range("J9").Select()
